# Threat Alert Report update (2026-01-18 01:00)
# The source data table (A2:K24) is refreshed with a new pull of threat
# comparisons; one obsolete row is removed so the table now runs A2:K23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "18-JAN-26"
$ws.Range("B2").Value = "SM-441"
$ws.Range("C2").Value = "Nesma Airlines NE-180"
$ws.Range("D2").Value = 8062
$ws.Range("E2").Value = 8087
$ws.Range("F2").Value = -25
$ws.Range("G2").Value = 30
$ws.Range("H2").Value = 30
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = "LOW THREAT"
$ws.Range("K2").Value = "EGP"
# Row 3
$ws.Range("A3").Value = "27-JAN-26"
$ws.Range("B3").Value = "SM-441"
$ws.Range("C3").Value = "Nile Air NP-103"
$ws.Range("D3").Value = 8742
$ws.Range("E3").Value = 9170
$ws.Range("F3").Value = -428
$ws.Range("G3").Value = 30
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = "LOW THREAT"
$ws.Range("K3").Value = "EGP"
# Row 4
$ws.Range("A4").Value = "29-JAN-26"
$ws.Range("B4").Value = "SM-441"
$ws.Range("C4").Value = "Nile Air NP-103"
$ws.Range("D4").Value = 8137
$ws.Range("E4").Value = 8628
$ws.Range("F4").Value = -491
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = "LOW THREAT"
$ws.Range("K4").Value = "EGP"
# Row 5
$ws.Range("A5").Value = "03-FEB-26"
$ws.Range("B5").Value = "SM-441"
$ws.Range("C5").Value = "Nesma Airlines NE-180"
$ws.Range("D5").Value = 7155
$ws.Range("E5").Value = 8628
$ws.Range("F5").Value = -1473
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 30
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "LOW THREAT"
$ws.Range("K5").Value = "EGP"
# Row 6
$ws.Range("A6").Value = "03-FEB-26"
$ws.Range("B6").Value = "SM-441"
$ws.Range("C6").Value = "Nile Air NP-103"
$ws.Range("D6").Value = 7658
$ws.Range("E6").Value = 8628
$ws.Range("F6").Value = -970
$ws.Range("G6").Value = 30
$ws.Range("H6").Value = 30
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "LOW THREAT"
$ws.Range("K6").Value = "EGP"
# Row 7
$ws.Range("A7").Value = "04-FEB-26"
$ws.Range("B7").Value = "SM-441"
$ws.Range("C7").Value = "Nesma Airlines NE-180"
$ws.Range("D7").Value = 6651
$ws.Range("E7").Value = 9170
$ws.Range("F7").Value = -2519
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K7").Value = "EGP"
# Row 8
$ws.Range("A8").Value = "05-FEB-26"
$ws.Range("B8").Value = "SM-441"
$ws.Range("C8").Value = "Nile Air NP-103"
$ws.Range("D8").Value = 6701
$ws.Range("E8").Value = 7545
$ws.Range("F8").Value = -844
$ws.Range("G8").Value = 30
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = "LOW THREAT"
$ws.Range("K8").Value = "EGP"
# Row 9
$ws.Range("A9").Value = "06-FEB-26"
$ws.Range("B9").Value = "SM-441"
$ws.Range("C9").Value = "Nile Air NP-103"
$ws.Range("D9").Value = 7658
$ws.Range("E9").Value = 8628
$ws.Range("F9").Value = -970
$ws.Range("G9").Value = 30
$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "LOW THREAT"
$ws.Range("K9").Value = "EGP"
# Row 10
$ws.Range("A10").Value = "07-FEB-26"
$ws.Range("B10").Value = "SM-441"
$ws.Range("C10").Value = "Nile Air NP-303"
$ws.Range("D10").Value = 7658
$ws.Range("E10").Value = 8628
$ws.Range("F10").Value = -970
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = 30
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = "LOW THREAT"
$ws.Range("K10").Value = "EGP"
# Row 11
$ws.Range("A11").Value = "07-FEB-26"
$ws.Range("B11").Value = "SM-441"
$ws.Range("C11").Value = "Nile Air NP-103"
$ws.Range("D11").Value = 7658
$ws.Range("E11").Value = 8628
$ws.Range("F11").Value = -970
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = 30
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = "LOW THREAT"
$ws.Range("K11").Value = "EGP"
# Row 12
$ws.Range("A12").Value = "22-FEB-26"
$ws.Range("B12").Value = "SM-987"
$ws.Range("C12").Value = "Nile Air NP-303"
$ws.Range("D12").Value = 17723
$ws.Range("E12").Value = 20621
$ws.Range("F12").Value = -2898
$ws.Range("G12").Value = 30
$ws.Range("H12").Value = 30
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K12").Value = "EGP"
# Row 13
$ws.Range("A13").Value = "22-FEB-26"
$ws.Range("B13").Value = "SM-987"
$ws.Range("C13").Value = "Nile Air NP-103"
$ws.Range("D13").Value = 17723
$ws.Range("E13").Value = 20621
$ws.Range("F13").Value = -2898
$ws.Range("G13").Value = 30
$ws.Range("H13").Value = 30
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = "MEDIUM THREAT - MONITOR"
$ws.Range("K13").Value = "EGP"
# Row 14
$ws.Range("A14").Value = "19-MAR-26"
$ws.Range("B14").Value = "SM-987"
$ws.Range("C14").Value = "Nile Air NP-103"
$ws.Range("D14").Value = 8137
$ws.Range("E14").Value = 20621
$ws.Range("F14").Value = -12484
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = 30
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K14").Value = "EGP"
# Row 15
$ws.Range("A15").Value = "24-MAR-26"
$ws.Range("B15").Value = "SM-441"
$ws.Range("C15").Value = "Nile Air NP-303"
$ws.Range("D15").Value = 8742
$ws.Range("E15").Value = 9170
$ws.Range("F15").Value = -428
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = 30
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = "LOW THREAT"
$ws.Range("K15").Value = "EGP"
# Row 16
$ws.Range("A16").Value = "24-MAR-26"
$ws.Range("B16").Value = "SM-441"
$ws.Range("C16").Value = "Nile Air NP-403"
$ws.Range("D16").Value = 8742
$ws.Range("E16").Value = 9170
$ws.Range("F16").Value = -428
$ws.Range("G16").Value = 30
$ws.Range("H16").Value = 30
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = "LOW THREAT"
$ws.Range("K16").Value = "EGP"
# Row 17
$ws.Range("A17").Value = "24-MAR-26"
$ws.Range("B17").Value = "SM-441"
$ws.Range("C17").Value = "Nile Air NP-103"
$ws.Range("D17").Value = 8742
$ws.Range("E17").Value = 9170
$ws.Range("F17").Value = -428
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = 30
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = "LOW THREAT"
$ws.Range("K17").Value = "EGP"
# Row 18
$ws.Range("A18").Value = "24-MAR-26"
$ws.Range("B18").Value = "SM-943"
$ws.Range("C18").Value = "Nile Air NP-303"
$ws.Range("D18").Value = 8742
$ws.Range("E18").Value = 20621
$ws.Range("F18").Value = -11879
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 30
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K18").Value = "EGP"
# Row 19
$ws.Range("A19").Value = "24-MAR-26"
$ws.Range("B19").Value = "SM-943"
$ws.Range("C19").Value = "Nile Air NP-403"
$ws.Range("D19").Value = 8742
$ws.Range("E19").Value = 20621
$ws.Range("F19").Value = -11879
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = 30
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K19").Value = "EGP"
# Row 20
$ws.Range("A20").Value = "24-MAR-26"
$ws.Range("B20").Value = "SM-943"
$ws.Range("C20").Value = "Nile Air NP-103"
$ws.Range("D20").Value = 8742
$ws.Range("E20").Value = 20621
$ws.Range("F20").Value = -11879
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 30
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K20").Value = "EGP"
# Row 21
$ws.Range("A21").Value = "25-MAR-26"
$ws.Range("B21").Value = "SM-441"
$ws.Range("C21").Value = "Nile Air NP-113"
$ws.Range("D21").Value = 8742
$ws.Range("E21").Value = 9762
$ws.Range("F21").Value = -1020
$ws.Range("G21").Value = 30
$ws.Range("H21").Value = 30
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = "LOW THREAT"
$ws.Range("K21").Value = "EGP"
# Row 22
$ws.Range("A22").Value = "25-MAR-26"
$ws.Range("B22").Value = "SM-987"
$ws.Range("C22").Value = "Nile Air NP-113"
$ws.Range("D22").Value = 8742
$ws.Range("E22").Value = 20621
$ws.Range("F22").Value = -11879
$ws.Range("G22").Value = 30
$ws.Range("H22").Value = 30
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K22").Value = "EGP"
# Row 23
$ws.Range("A23").Value = "26-MAR-26"
$ws.Range("B23").Value = "SM-987"
$ws.Range("C23").Value = "Nile Air NP-103"
$ws.Range("D23").Value = 11979
$ws.Range("E23").Value = 20621
$ws.Range("F23").Value = -8642
$ws.Range("G23").Value = 30
$ws.Range("H23").Value = 30
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = "HIGH THREAT ALERT - NEED ACTION"
$ws.Range("K23").Value = "EGP"

# --- Fix up the IMPACT (column J) conditional-style highlighting ---
# J3 / J13 / J19 keep the same LOW/MEDIUM/HIGH threat level both
# before and after this refresh, so they still carry the correct
# untouched cell style (s=3/4/5 respectively). Copy just the
# formatting (not the value) from each onto every J cell whose
# threat level changed with the refreshed figures.
$ws.Range("J3").Copy()
$ws.Range("J8").PasteSpecial(-4122)
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J21").PasteSpecial(-4122)

$ws.Range("J13").Copy()
$ws.Range("J7").PasteSpecial(-4122)
$ws.Range("J12").PasteSpecial(-4122)

$ws.Range("J19").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("J22").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Drop the now-obsolete last row (old row 24) ---
$ws.Rows.Item(24).Delete()
